$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 99
$ws.Range("H99").Value = 326
$ws.Range("I99").Value = 326
$ws.Range("K99").Value = 978
$ws.Range("M99").Value = 520
# Row 115
$ws.Range("H115").Value = 5161.8335
$ws.Range("I115").Value = 5161.8335
$ws.Range("K115").Value = 15485.5005
$ws.Range("M115").Value = -13918.5005
# Row 116
$ws.Range("H116").Value = 4352.5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4352.5
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 4352.5
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -11236.5
# Row 118
$ws.Range("H118").Value = 1825
$ws.Range("I118").Value = 1500
$ws.Range("K118").Value = 4500
$ws.Range("M118").Value = -2843
# Row 125
$ws.Range("H125").Value = 1972
$ws.Range("I125").Value = 1972
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 17748
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -15288
$ws.Range("N125").ClearContents()
# Row 132
$ws.Range("H132").Value = 1704.875
$ws.Range("I132").Value = 1656.9166
$ws.Range("J132").Value = 1848.75
$ws.Range("K132").Value = 4970.7498
$ws.Range("L132").Value = 5546.25
$ws.Range("M132").Value = -2440.7498
$ws.Range("N132").Value = -10606.25
# Row 137
$ws.Range("H137").Value = 3941.5557
$ws.Range("I137").Value = 3809.375
$ws.Range("K137").Value = 11428.125
$ws.Range("M137").Value = -8878.125
# Row 138
$ws.Range("H138").Value = 5130.552
$ws.Range("I138").Value = 2683.1428
$ws.Range("J138").Value = 5909.273
$ws.Range("K138").Value = 8049.428400000001
$ws.Range("L138").Value = 17727.819
$ws.Range("M138").Value = -2909.428400000001
$ws.Range("N138").Value = -28007.819

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1019.13336
$ws.Range("I61").Value = 877.6429000000001
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 877.6429000000001
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -665.6429000000001
$ws.Range("N61").Value = -3424
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
# Row 74
$ws.Range("H74").Value = 5128932.5
$ws.Range("I74").Value = 6451370
$ws.Range("K74").Value = 6451370
$ws.Range("M74").Value = -6450496
# Row 77
$ws.Range("H77").Value = 5128932.5
$ws.Range("I77").Value = 6451370
$ws.Range("K77").Value = 32256850
$ws.Range("M77").Value = -32252482
# Row 122
$ws.Range("H122").Value = 3599.2
$ws.Range("I122").Value = 3599.2
$ws.Range("K122").Value = 10797.6
$ws.Range("M122").Value = -8347.599999999999
# Row 132
$ws.Range("H132").Value = 1794.5193
$ws.Range("I132").Value = 863.7895
$ws.Range("J132").Value = 4320.7856
$ws.Range("K132").Value = 2591.3685
$ws.Range("L132").Value = 12962.3568
$ws.Range("M132").Value = -61.36850000000004
$ws.Range("N132").Value = -18022.3568
# Row 136
$ws.Range("H136").Value = 1019.13336
$ws.Range("I136").Value = 877.6429000000001
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 2632.9287
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -82.92870000000039
$ws.Range("N136").Value = -14100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 4990
$ws.Range("I105").Value = 4990
$ws.Range("K105").Value = 4990
$ws.Range("M105").Value = -3243
# Row 134
$ws.Range("H134").Value = 3370.4583
$ws.Range("I134").Value = 3370.4583
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10111.3749
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7576.374899999999
$ws.Range("N134").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2442.9744
$ws.Range("I31").Value = 2340.6667
$ws.Range("J31").Value = 3670.6667
$ws.Range("K31").Value = 2340.6667
$ws.Range("L31").Value = 3670.6667
$ws.Range("M31").Value = -2045.6667
$ws.Range("N31").Value = -4260.6667
# Row 34
$ws.Range("H34").Value = 2442.9744
$ws.Range("I34").Value = 2340.6667
$ws.Range("J34").Value = 3670.6667
$ws.Range("K34").Value = 2340.6667
$ws.Range("L34").Value = 3670.6667
$ws.Range("M34").Value = -2138.6667
$ws.Range("N34").Value = -4074.6667
# Row 58
$ws.Range("H58").Value = 1644.9286
$ws.Range("I58").Value = 1271.6364
$ws.Range("K58").Value = 1271.6364
$ws.Range("M58").Value = -1068.6364
# Row 134
$ws.Range("H134").Value = 2174.2
$ws.Range("I134").Value = 2236.6428
$ws.Range("J134").Value = 1300
$ws.Range("K134").Value = 6709.928400000001
$ws.Range("L134").Value = 3900
$ws.Range("M134").Value = -4174.928400000001
$ws.Range("N134").Value = -8970
# Row 136
$ws.Range("H136").Value = 1644.9286
$ws.Range("I136").Value = 1271.6364
$ws.Range("K136").Value = 3814.9092
$ws.Range("M136").Value = -1264.9092

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 497.4
$ws.Range("I8").Value = 497.4
$ws.Range("K8").Value = 1492.2
$ws.Range("M8").Value = -1353.2
# Row 113
$ws.Range("H113").Value = 1084.85
$ws.Range("I113").Value = 1224.4
$ws.Range("J113").Value = 1038.3334
$ws.Range("K113").Value = 3673.2
$ws.Range("L113").Value = 3115.0002
$ws.Range("M113").Value = -1503.2
$ws.Range("N113").Value = -7455.0002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 30000
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
# Row 37
$ws.Range("H37").Value = 30000
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
# Row 102
$ws.Range("H102").Value = 2434.3333
$ws.Range("I102").Value = 842.1429000000001
$ws.Range("J102").Value = 8007
$ws.Range("K102").Value = 842.1429000000001
$ws.Range("L102").Value = 8007
$ws.Range("M102").Value = 779.8570999999999
$ws.Range("N102").Value = -11251
# Row 113
$ws.Range("H113").Value = 2490
$ws.Range("I113").Value = 2490
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2490
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -320
$ws.Range("N113").ClearContents()
# Row 132
$ws.Range("H132").Value = 1348.0385
$ws.Range("I132").Value = 871.8261
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 2615.4783
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -85.47829999999976
$ws.Range("N132").Value = -20057

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3015.9
$ws.Range("I7").Value = 2962.111
$ws.Range("K7").Value = 2962.111
$ws.Range("M7").Value = -2850.111
# Row 22
$ws.Range("H22").Value = 753.3333
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 753.3333
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 753.3333
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1343.3333
# Row 27
$ws.Range("H27").Value = 753.3333
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 753.3333
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 753.3333
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -967.3333
# Row 126
$ws.Range("H126").Value = 3015.9
$ws.Range("I126").Value = 2962.111
$ws.Range("K126").Value = 8886.332999999999
$ws.Range("M126").Value = -6416.332999999999
# Row 136
$ws.Range("H136").Value = 10002912
$ws.Range("I136").Value = 11431042
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 34293126
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -34290576
$ws.Range("N136").Value = -23100
